# Update "Rally" take-rate from 26% to 25%, and remove two obsolete
# "Optional" rows (RS / Optional Bundle 16, RS / Optional Bundle 7)
# from the db_take_rate table on Foglio1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")
$ws.Activate()

# --- 1) Rally take-rate 0.26 -> 0.25 (column B, rows where column A = "Rally") ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row()
for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    if ($name -eq "Rally") {
        $ws.Cells.Item($r, 2).Value = 0.25
    }
}

# --- 2) Delete the "RS" / "Optional Bundle 7" row (further down, deleted first) ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row()
for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    $opt = $ws.Cells.Item($r, 3).Value()
    if ($name -eq "RS" -and $opt -eq "Optional Bundle 7") {
        $row = $ws.Rows.Item($r)
        $row.Select()
        $row.Delete()
        break
    }
}

# --- 3) Delete the "RS" / "Optional Bundle 16" row (deleted last, leaves selection here) ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row()
for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    $opt = $ws.Cells.Item($r, 3).Value()
    if ($name -eq "RS" -and $opt -eq "Optional Bundle 16") {
        $row = $ws.Rows.Item($r)
        $row.Select()
        $row.Delete()
        break
    }
}

$wb.Save()
